# Apply updated 焦煤上游总库存_预测 (column C) and 焦煤上游总库存 (column B) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (forecast) values newly populated / revised for rows 7-21
$ws.Range("C7").Value  = 960.8
$ws.Range("C8").Value  = 969.3
$ws.Range("C9").Value  = 969
$ws.Range("C10").Value = 979.4
$ws.Range("C11").Value = 977.8
$ws.Range("C12").Value = 960.6
$ws.Range("C13").Value = 950.7
$ws.Range("C14").Value = 947.6
$ws.Range("C15").Value = 948.7
$ws.Range("C16").Value = 955.9
$ws.Range("C17").Value = 939.6
$ws.Range("C18").Value = 950.6
$ws.Range("C19").Value = 953
$ws.Range("C20").Value = 941.6
$ws.Range("C21").Value = 940.8

# Column B (actual) values revised for rows 22-23
$ws.Range("B22").Value = 870.642
$ws.Range("B23").Value = 890.884
